$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to reflect the latest scrape.
# Price-column values are forced to Text format so Excel keeps trailing zeros
# and multi-dot "thousands" groupings (e.g. "1.828.39") exactly as scraped.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.738.42"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.602.32"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.82"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.72"
$ws.Range("E10").Value = "  +0.99%  "

$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.39"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.600.11"
$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.12"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0740"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "210.23"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("E20").Value = "  +2.67%  "

$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("E22").Value = "  -2.58%  "

$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.66"
$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.10"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.39"
$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0511"
$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.294.14"
$ws.Range("E33").Value = "  +0.74%  "

$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.605"
$ws.Range("E36").Value = "  -2.11%  "

$ws.Range("E37").Value = "  +10.94%  "

$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("E40").Value = "  -2.11%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.82"
$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.739.38"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.74"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.55"
$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("E48").Value = "  +1.72%  "

$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.44"
$ws.Range("E50").Value = "  +0.90%  "

$ws.Range("E51").Value = "  +0.91%  "
